# Generate Report for Handback
# - Updates the existing "9db4f506..." row data to the "eb4d2096..." file
#   (new hash/timestamps for the zh-cn + de-de handback), and
# - Appends a new row for the "fe0a9c1a..." file on all three sheets
#   (Overview, zh-cn, de-de), including hyperlinks + table range growth.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet (A1:G2 -> A1:G3)
# ---------------------------------------------------------------------

# Row 2 (existing file) gets the renamed uuid + refreshed HO generate date
$wsOverview.Range("A2").Value = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.md"
$wsOverview.Range("B2").Value = "e2e\eb4d2096-cf7a-4f78-a4ec-d67885d88521.md"
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G2").Value = "2016-08-28 11:09:35"

# Row 3 (new file)
$wsOverview.Range("A3").Value = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md"
$wsOverview.Range("B3").Value = "e2e\fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-28 11:09:35"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f364f6003c6ee15ff7d8e3199207c0783deec7/e2e/fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md", $null, $null, "e2e\fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md") | Out-Null
$wsOverview.Range("B3").Font.Underline = 2
$wsOverview.Range("B3").Font.Color = 15570276

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# zh-cn sheet (A1:P2 -> A1:P3)
# ---------------------------------------------------------------------

# Row 2 (existing file) refreshed hash / timestamps
$wsZhCn.Range("A2").Value = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D2").Value = "e2e"
$wsZhCn.Range("E2").Value = "ht"
$wsZhCn.Range("F2").Value = "False"
$wsZhCn.Range("G2").Value = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.5387253624e2db7618f5b0610dc21bb69a71115a.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-28 11:09:30"
$wsZhCn.Range("I2").Value = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.md"
$wsZhCn.Range("J2").Value = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.5387253624e2db7618f5b0610dc21bb69a71115a.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-28 11:09:47"
$wsZhCn.Range("L2").Value = ""
$wsZhCn.Range("M2").Value = "True"
$wsZhCn.Range("N2").Value = ""
$wsZhCn.Range("O2").Value = "False"
$wsZhCn.Range("P2").Value = ""

# Row 3 (new file)
$wsZhCn.Range("A3").Value = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.8271afb11851a3eb9dd5f5ed7a864122870a8f1d.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-28 11:09:30"
$wsZhCn.Range("I3").Value = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md"
$wsZhCn.Range("J3").Value = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.8271afb11851a3eb9dd5f5ed7a864122870a8f1d.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-28 11:09:47"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f364f6003c6ee15ff7d8e3199207c0783deec7/e2e/fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md", $null, $null, "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md") | Out-Null
$wsZhCn.Range("A3").Font.Underline = 2
$wsZhCn.Range("A3").Font.Color = 15570276

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d9bbc94ef822cf502d57c352d1fb903dcfafde2c/e2e/fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md", $null, $null, "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md") | Out-Null
$wsZhCn.Range("I3").Font.Underline = 2
$wsZhCn.Range("I3").Font.Color = 15570276

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------
# de-de sheet (A1:P2 -> A1:P3)
# ---------------------------------------------------------------------

# Row 2 (existing file) refreshed hash / timestamps
$wsDeDe.Range("A2").Value = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D2").Value = "e2e"
$wsDeDe.Range("E2").Value = "ht"
$wsDeDe.Range("F2").Value = "False"
$wsDeDe.Range("G2").Value = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.5387253624e2db7618f5b0610dc21bb69a71115a.de-de.xlf"
$wsDeDe.Range("I2").Value = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.md"
$wsDeDe.Range("J2").Value = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.5387253624e2db7618f5b0610dc21bb69a71115a.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-28 11:09:54"
$wsDeDe.Range("L2").Value = ""
$wsDeDe.Range("M2").Value = "True"
$wsDeDe.Range("N2").Value = ""
$wsDeDe.Range("O2").Value = "False"
$wsDeDe.Range("P2").Value = ""

# Row 3 (new file)
$wsDeDe.Range("A3").Value = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.8271afb11851a3eb9dd5f5ed7a864122870a8f1d.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-28 11:09:35"
$wsDeDe.Range("I3").Value = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md"
$wsDeDe.Range("J3").Value = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.8271afb11851a3eb9dd5f5ed7a864122870a8f1d.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-28 11:09:54"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f364f6003c6ee15ff7d8e3199207c0783deec7/e2e/fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md", $null, $null, "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md") | Out-Null
$wsDeDe.Range("A3").Font.Underline = 2
$wsDeDe.Range("A3").Font.Color = 15570276

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/efb43ddee6d7c1f5ffdd0d4e01baf1345331e9ce/e2e/fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md", $null, $null, "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md") | Out-Null
$wsDeDe.Range("I3").Font.Underline = 2
$wsDeDe.Range("I3").Font.Color = 15570276

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

Write-Output "Handback report rows updated."
